$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.977.90'
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = '2.548.38'
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.26'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.41'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.05%  '

$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -0.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.65'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0816'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.70%  '

$ws.Range("E12").Value = '  -0.20%  '

$ws.Range("E13").Value = '  +5.85%  '

$ws.Range("D14").Value = '2.941.97'
$ws.Range("E14").Value = '  -0.25%  '

$ws.Range("D15").Value = '2.557.77'
$ws.Range("E15").Value = '  +0.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.879'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.47%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.80'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +4.20%  '

$ws.Range("D18").Value = '43.169.83'
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.62'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.61%  '

$ws.Range("D20").Value = '0.0₃0987'
$ws.Range("E20").Value = '  +0.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.60'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.87'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.74%  '

$ws.Range("E24").Value = '  +1.45%  '

$ws.Range("E25").Value = '  -2.10%  '

$ws.Range("E26").Value = '  -5.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("E28").Value = '  +1.21%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.90'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.89%  '

$ws.Range("E30").Value = '  -1.17%  '

$ws.Range("E31").Value = '  +0.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '159.28'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.09%  '

$ws.Range("E33").Value = '  -0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.16'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("E35").Value = '  +0.83%  '

$ws.Range("E36").Value = '  -2.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.87'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +12.63%  '

$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.75'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +10.36%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.115'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.95%  '

$ws.Range("E40").Value = '  -0.56%  '

$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.43'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.63%  '

$ws.Range("B42").Value = 'ApeXProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.08'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +30.99%  '

$ws.Range("E43").Value = '  -0.26%  '

$ws.Range("D44").Value = '2.102.14'
$ws.Range("E44").Value = '  +1.65%  '

$ws.Range("E45").Value = '  -2.39%  '

$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '86.37'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.59%  '

$ws.Range("E48").Value = '  +2.86%  '

$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.799.54'
$ws.Range("E49").Value = '  -0.21%  '

$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.14'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +8.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '103.64'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.69%  '
